$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# --- Crime Complaints table updates (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 50
$ws.Range("F14").Value = 10
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = 25
$ws.Range("I14").Value = 89
$ws.Range("J14").Value = 97
$ws.Range("K14").Value = -8.247422680412
$ws.Range("L14").Value = -14.423076923076
$ws.Range("M14").Value = -5.31914893617
$ws.Range("N14").Value = -73.511904761904

# Row 15
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 33.333333333333
$ws.Range("F15").Value = 33
$ws.Range("G15").Value = 25
$ws.Range("H15").Value = 32
$ws.Range("I15").Value = 263
$ws.Range("J15").Value = 271
$ws.Range("K15").Value = -2.952029520295
$ws.Range("L15").Value = 9.128630705394
$ws.Range("M15").Value = 30.19801980198
$ws.Range("N15").Value = -46.975806451612

# Row 16
$ws.Range("C16").Value = 102
$ws.Range("D16").Value = 106
$ws.Range("E16").Value = -3.77358490566
$ws.Range("F16").Value = 436
$ws.Range("G16").Value = 433
$ws.Range("H16").Value = 0.692840646651
$ws.Range("I16").Value = 3269
$ws.Range("J16").Value = 3464
$ws.Range("K16").Value = -5.629330254041
$ws.Range("L16").Value = 34.915394139496
$ws.Range("M16").Value = 12.33676975945
$ws.Range("N16").Value = -69.627427297222

# Row 17
$ws.Range("C17").Value = 174
$ws.Range("D17").Value = 153
$ws.Range("E17").Value = 13.725490196078
$ws.Range("F17").Value = 673
$ws.Range("G17").Value = 583
$ws.Range("H17").Value = 15.437392795883
$ws.Range("I17").Value = 5471
$ws.Range("J17").Value = 4988
$ws.Range("K17").Value = 9.683239775461
$ws.Range("L17").Value = 33.896231032794
$ws.Range("M17").Value = 81.278992710404
$ws.Range("N17").Value = -12.67358339984

# Row 18
$ws.Range("C18").Value = 49
$ws.Range("D18").Value = 64
$ws.Range("E18").Value = -23.4375
$ws.Range("F18").Value = 214
$ws.Range("G18").Value = 216
$ws.Range("H18").Value = -0.925925925925
$ws.Range("I18").Value = 2011
$ws.Range("J18").Value = 1974
$ws.Range("K18").Value = 1.874366767983
$ws.Range("L18").Value = 41.420534458509
$ws.Range("M18").Value = -6.638811513463
$ws.Range("N18").Value = -84.181546448517

# Row 19
$ws.Range("C19").Value = 146
$ws.Range("D19").Value = 167
$ws.Range("E19").Value = -12.574850299401
$ws.Range("F19").Value = 689
$ws.Range("G19").Value = 639
$ws.Range("H19").Value = 7.824726134585
$ws.Range("I19").Value = 5322
$ws.Range("J19").Value = 5378
$ws.Range("K19").Value = -1.041279285979
$ws.Range("L19").Value = 25.017618040873
$ws.Range("M19").Value = 70.905587668593
$ws.Range("N19").Value = 4.434850863422

# Row 20
$ws.Range("C20").Value = 97
$ws.Range("D20").Value = 62
$ws.Range("E20").Value = 56.451612903225
$ws.Range("F20").Value = 406
$ws.Range("G20").Value = 262
$ws.Range("H20").Value = 54.961832061068
$ws.Range("I20").Value = 3609
$ws.Range("J20").Value = 2631
$ws.Range("K20").Value = 37.172177879133
$ws.Range("L20").Value = 96.675749318801
$ws.Range("M20").Value = 156.503198294243
$ws.Range("N20").Value = -64.998545242944

# Row 21
$ws.Range("C21").Value = 579
$ws.Range("D21").Value = 560
$ws.Range("E21").Value = 3.392857142857
$ws.Range("F21").Value = 2461
$ws.Range("G21").Value = 2166
$ws.Range("H21").Value = 13.619575253924
$ws.Range("I21").Value = 20034
$ws.Range("J21").Value = 18803
$ws.Range("K21").Value = 6.546827633888
$ws.Range("L21").Value = 39.434855233853
$ws.Range("M21").Value = 55.314365454686
$ws.Range("N21").Value = -56.428882122662

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 9
$ws.Range("E22").Value = -44.444444444444
$ws.Range("F22").Value = 21
$ws.Range("G22").Value = 32
$ws.Range("H22").Value = -34.375
$ws.Range("I22").Value = 192
$ws.Range("J22").Value = 244
$ws.Range("K22").Value = -21.311475409836
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = -9.43396226415

# Row 23
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 30
$ws.Range("E23").Value = 6.666666666666
$ws.Range("F23").Value = 149
$ws.Range("G23").Value = 121
$ws.Range("H23").Value = 23.140495867768
$ws.Range("I23").Value = 1215
$ws.Range("J23").Value = 1080
$ws.Range("K23").Value = 12.5
$ws.Range("L23").Value = 52.255639097744
$ws.Range("M23").Value = 67.125171939477

# Row 24
$ws.Range("C24").Value = 353
$ws.Range("D24").Value = 370
$ws.Range("E24").Value = -4.594594594594
$ws.Range("F24").Value = 1425
$ws.Range("G24").Value = 1510
$ws.Range("H24").Value = -5.629139072847
$ws.Range("I24").Value = 12101
$ws.Range("J24").Value = 12514
$ws.Range("K24").Value = -3.3003036599
$ws.Range("L24").Value = 42.987120406475
$ws.Range("M24").Value = 40.905915230554

# Row 25
$ws.Range("C25").Value = 222
$ws.Range("D25").Value = 178
$ws.Range("E25").Value = 24.719101123595
$ws.Range("F25").Value = 831
$ws.Range("G25").Value = 721
$ws.Range("H25").Value = 15.256588072122
$ws.Range("I25").Value = 7159
$ws.Range("J25").Value = 6801
$ws.Range("K25").Value = 5.263931774739
$ws.Range("L25").Value = 26.283295113776
$ws.Range("M25").Value = -5.090812674002

# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 56
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = 21.739130434782
$ws.Range("I26").Value = 443
$ws.Range("J26").Value = 468
$ws.Range("K26").Value = -5.34188034188
$ws.Range("L26").Value = 12.151898734177

# Row 27
$ws.Range("C27").Value = 22
$ws.Range("D27").Value = 15
$ws.Range("E27").Value = 46.666666666666
$ws.Range("F27").Value = 74
$ws.Range("G27").Value = 78
$ws.Range("H27").Value = -5.128205128205
$ws.Range("I27").Value = 698
$ws.Range("J27").Value = 615
$ws.Range("K27").Value = 13.495934959349
$ws.Range("L27").Value = 15.181518151815

# Row 28
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 40
$ws.Range("F28").Value = 35
$ws.Range("G28").Value = 33
$ws.Range("H28").Value = 6.060606060606
$ws.Range("I28").Value = 279
$ws.Range("J28").Value = 352
$ws.Range("K28").Value = -20.738636363636
$ws.Range("L28").Value = -31.784841075794
$ws.Range("M28").Value = -16.964285714285
$ws.Range("N28").Value = -70.600632244467

# Row 29
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 20
$ws.Range("F29").Value = 31
$ws.Range("G29").Value = 27
$ws.Range("H29").Value = 14.814814814814
$ws.Range("I29").Value = 229
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = -23.666666666666
$ws.Range("L29").Value = -34.383954154727
$ws.Range("M29").Value = -18.214285714285
$ws.Range("N29").Value = -73.49537037037

# Row 30
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -66.666666666666
$ws.Range("L30").Value = -57.575757575757

